$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 125, shifting existing rows 125-252 down to 126-253.
$ws.Rows(125).Insert()

# Populate the newly inserted row 125 with the new weekly record.
# (same market/product metadata as the record that used to sit here,
#  but a new date and updated price figures)
$ws.Cells.Item(125, 1).Value = 3
$ws.Cells.Item(125, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(125, 3).Value = "Coquimbo"
$ws.Cells.Item(125, 4).Value = "2021-11-12"
$ws.Cells.Item(125, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(125, 5).Value = 5
$ws.Cells.Item(125, 6).Value = "Fruta"
$ws.Cells.Item(125, 7).Value = 100108
$ws.Cells.Item(125, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(125, 9).Value = 100108002
$ws.Cells.Item(125, 10).Value = "Mango"
$ws.Cells.Item(125, 11).Value = "Sin especificar"
$ws.Cells.Item(125, 12).Value = "Primera"
$ws.Cells.Item(125, 13).Value = 456
$ws.Cells.Item(125, 14).Value = 6000
$ws.Cells.Item(125, 15).Value = 6000
$ws.Cells.Item(125, 16).Value = 6000
$ws.Cells.Item(125, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(125, 18).Value = "Perú"
$ws.Cells.Item(125, 19).Value = 1500
$ws.Cells.Item(125, 20).Value = 4
